# Updates the "cryptos" price list with refreshed Price / Volume(1h) figures.
# (Row 43/44 also swap Coin/Link: EnergySwap and Monero trade places.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.712.60'
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").Value = '2.783.85'
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''352.77'
$ws.Range("E5").Value = '  -1.53%  '
$ws.Range("D6").Value = '''109.02'
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("E7").Value = '  -2.25%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("D9").Value = '''0.604'
$ws.Range("E9").Value = '  +1.97%  '
$ws.Range("D10").Value = '''39.88'
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("E11").Value = '  +2.57%  '
$ws.Range("D12").Value = '''20.20'
$ws.Range("E12").Value = '  +3.63%  '
$ws.Range("E13").Value = '  -2.00%  '
$ws.Range("E14").Value = '  +1.35%  '
$ws.Range("D15").Value = '3.222.99'
$ws.Range("E15").Value = '  -0.13%  '
$ws.Range("D16").Value = '2.789.26'
$ws.Range("E16").Value = '  +0.46%  '
$ws.Range("E17").Value = '  -2.38%  '
$ws.Range("D18").Value = '51.736.10'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("D19").Value = '''7.75'
$ws.Range("E19").Value = '  +4.57%  '
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").Value = '''13.16'
$ws.Range("E21").Value = '  +1.27%  '
$ws.Range("D22").Value = '0.0₃0966'
$ws.Range("E22").Value = '  -1.59%  '
$ws.Range("D23").Value = '''69.91'
$ws.Range("E23").Value = '  -0.50%  '
$ws.Range("D24").Value = '''266.93'
$ws.Range("E24").Value = '  -2.69%  '
$ws.Range("D25").Value = '''2.74'
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").Value = '''26.15'
$ws.Range("E26").Value = '  -2.00%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").Value = '''0.163'
$ws.Range("E28").Value = '  +12.26%  '
$ws.Range("D29").Value = '''10.22'
$ws.Range("E29").Value = '  +0.32%  '
$ws.Range("D30").Value = '''37.29'
$ws.Range("E30").Value = '  +7.72%  '
$ws.Range("E31").Value = '  -2.50%  '
$ws.Range("D32").Value = '''6.19'
$ws.Range("E32").Value = '  +8.14%  '
$ws.Range("D33").Value = '''51.75'
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("D34").Value = '''0.0454'
$ws.Range("E34").Value = '  -2.29%  '
$ws.Range("D35").Value = '''5.55'
$ws.Range("E35").Value = '  +4.91%  '
$ws.Range("D36").Value = '''0.0832'
$ws.Range("E36").Value = '  -1.61%  '
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("D38").Value = '''18.51'
$ws.Range("E38").Value = '  +2.76%  '
$ws.Range("E39").Value = '  -2.52%  '
$ws.Range("E40").Value = '  -1.51%  '
$ws.Range("D41").Value = '''2.54'
$ws.Range("E41").Value = '  -1.12%  '
$ws.Range("E42").Value = '  -0.55%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '''22.19'
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").Value = '''120.45'
$ws.Range("E44").Value = '  -1.48%  '
$ws.Range("E45").Value = '  -2.83%  '
$ws.Range("D46").Value = '2.122.08'
$ws.Range("E46").Value = '  +2.23%  '
$ws.Range("E47").Value = '  +1.74%  '
$ws.Range("E48").Value = '  +6.84%  '
$ws.Range("D49").Value = '''0.911'
$ws.Range("E49").Value = '  -2.24%  '
$ws.Range("D50").Value = '''5.44'
$ws.Range("E50").Value = '  -4.91%  '
$ws.Range("E51").Value = '  +9.03%  '
